$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.613.99"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "3.491.71"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.74"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.88"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.490.74"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("E10").Value = "  +6.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.11"
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "4.098.52"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.08"
$ws.Range("E14").Value = "  +10.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.136"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").Value = "67.607.96"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "3.498.27"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.27"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.33"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.11"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.88"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.11"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.535"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.47"
$ws.Range("E28").Value = "  +1.91%  "
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("E31").Value = "  -1.80%  "
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.06"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.59"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.37"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.60"
$ws.Range("E37").Value = "  -3.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.82"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.873"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.91"
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("E42").Value = "  +6.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.65"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "2.850.57"
$ws.Range("E44").Value = "  +1.38%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.03"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0725"
$ws.Range("E46").Value = "  -2.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.49"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.05"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "335.63"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.05"
$ws.Range("E51").Value = "  -2.28%  "
